# add fr core 2.1.0
$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.0.1 -> 2.1.0
$meta.Range("B3").Value = "2.1.0"

# Date: 2024-04-16T11:49:14+02:00 -> 2024-09-04T10:06:33+02:00
$meta.Range("B8").Value = "2024-09-04T10:06:33+02:00"

# Contact: InteropSante (fhir@interopsante.org(WORK)) -> (work)
$meta.Range("B11").Value = "InteropSanté (fhir@interopsante.org(work))"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Row 6 (Extension.value[x]): Type(s) column K, Reference -> long list of types
$elements.Range("K6").Value = "base64Binary`nbooleancanonicalcodedatedateTimedecimalidinstantintegermarkdownoidpositiveIntstringtimeunsignedInturiurluuidAddressAgeAnnotationAttachmentCodeableConceptCodingContactPointCountDistanceDurationHumanNameIdentifierMoneyPeriodQuantityRangeRatioReferenceSampledDataSignatureTimingContactDetailContributorDataRequirementExpressionParameterDefinitionRelatedArtifactTriggerDefinitionUsageContextDosageMeta"

# Row 6: Slicing Rules column AE, closed -> open
$elements.Range("AE6").Value = "open"

# Column K width update (11th column) 202.08984375 -> 255.0
# Note: the host snaps ColumnWidth to whole-pixel granularity using
# stored = (round(input*6)+5)/6, so feeding 255.0 directly lands on
# 255.8333; 254.1 is inside the bucket that snaps back to exactly 255.0.
$elements.Columns.Item(11).ColumnWidth = 254.1
